$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price (D) and volume-change (E) values for the crypto list.
# Numeric-looking price strings get an apostrophe text-prefix (as Excel's UI
# would do for text entry) followed by a Style reset so no extra number-format
# style lingers on the cell - keeps the cell a plain text value like the source.

$ws.Range("D2").Value = "27.431.94"
$ws.Range("E2").Value = "  -3.13%  "
$ws.Range("D3").Value = "1.652.82"
$ws.Range("E3").Value = "  -3.75%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'214.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("D6").Value = "'0.511"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.10%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'24.13"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("D11").Value = "'0.0876"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("D12").Value = "1.887.50"
$ws.Range("E12").Value = "  -3.64%  "
$ws.Range("D13").Value = "1.656.47"
$ws.Range("E13").Value = "  -3.54%  "
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("D15").Value = "'0.564"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").Value = "'65.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.57%  "
$ws.Range("D17").Value = "27.448.15"
$ws.Range("E17").Value = "  -3.09%  "
$ws.Range("D18").Value = "'235.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.26%  "
$ws.Range("E19").Value = "  -2.73%  "
$ws.Range("D20").Value = "'7.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.08%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("E22").Value = "  -3.48%  "
$ws.Range("D23").Value = "'9.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.43%  "
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D25").Value = "'145.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("E26").Value = "  -2.97%  "
$ws.Range("D27").Value = "'16.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.63%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  -2.38%  "
$ws.Range("E30").Value = "  -2.77%  "
$ws.Range("D31").Value = "'1.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("E32").Value = "  -2.76%  "
$ws.Range("D33").Value = "1.449.11"
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("E34").Value = "  -4.41%  "
$ws.Range("E35").Value = "  -4.43%  "
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("D37").Value = "'0.914"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.90%  "
$ws.Range("E38").Value = "  -4.48%  "
$ws.Range("E39").Value = "  -2.92%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "'66.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.81%  "
$ws.Range("E43").Value = "  -3.52%  "
$ws.Range("E44").Value = "  -2.60%  "
$ws.Range("E45").Value = "  -1.84%  "
$ws.Range("D46").Value = "1.795.04"
$ws.Range("E46").Value = "  -3.60%  "
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "'88.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("D49").Value = "0.0₆0107"
$ws.Range("E49").Value = "  -6.71%  "
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("E51").Value = "  -3.23%  "
